$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B26 value
$ws.Range("B26").Value = 352

# Add new row 27
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 204

# Add new row 28
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 1
